# Rearrange module D to calculate emissions from activity and ef.
# activity and ef dbs are same size for nc and comb emissions.
#
# Adds a "type" column (D) to the Sectors sheet, classifying each sector
# row as combustion ("comb", units = kt) or non-combustion ("NC",
# units = B2005USD) based on the existing "units" column (C).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sectors")

# Header - match the style of the other header cells (C1).
$ws.Cells.Item(1, 3).Copy()
$ws.Cells.Item(1, 4).PasteSpecial(-4122)
$ws.Cells.Item(1, 4).Value = "type"

# Use the formatting already present on D2 (existing "comb" style) as the
# template for every data row, so odd rows (e.g. D37, which had drifted to
# a different style) get normalized back to the common look before we
# stamp in the actual value.
$ws.Cells.Item(2, 4).Copy()

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 59) { $lastRow = 59 }

for ($r = 2; $r -le $lastRow; $r++) {
    $target = $ws.Cells.Item($r, 4)
    $target.PasteSpecial(-4122)

    $units = $ws.Cells.Item($r, 3).Value2
    if ($units -eq "kt") {
        $target.Value = "comb"
    } else {
        $target.Value = "NC"
    }
}

$excel.CutCopyMode = 0

# Restore selection/view as left by the editor: active cell D5, no pinned
# top-left scroll position.
$ws.Activate() | Out-Null
$ws.Range("D5").Select() | Out-Null
